$d = $word.ActiveDocument

# --- Locate the "Q11: " run and split it into its own paragraph -----------
$rng = $d.Content
$rng.Find.Execute("Q11: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphAfter()

# --- Insert the first "Used this link..." note before the first hyperlink -
$pNote1 = $d.Paragraphs.Item(22)
$insertPoint1 = $d.Range($pNote1.Range.Start, $pNote1.Range.Start)
$insertPoint1.InsertBefore("Used this link to learn how to set up SDL2 and basic code to get a window appearing`r")

# --- Split off a new paragraph after the first hyperlink paragraph --------
$pLink1 = $d.Paragraphs.Item(23)
$pLink1.Range.InsertParagraphAfter()

# --- Insert the second "Used this link..." note ----------------------------
$pNote2 = $d.Paragraphs.Item(24)
$pNote2.Range.InsertBefore("Used this link to learn how to handle text input and keyboard input")

# --- Create a new empty paragraph to host the second hyperlink -------------
$pNote2b = $d.Paragraphs.Item(24)
$pNote2b.Range.InsertParagraphAfter()

# --- Add the second hyperlink (text input / keyboard input video) ---------
$pLink2 = $d.Paragraphs.Item(25)
$target2 = "https://www.youtube.com/watch?v=m2doh3Li65c&ab_channel=CodingMadeEasy"
$pLink2.Range.InsertBefore($target2)

$pLink2b = $d.Paragraphs.Item(25)
$linkRange2 = $d.Range($pLink2b.Range.Start, $pLink2b.Range.Start + $target2.Length)
$d.Hyperlinks.Add($linkRange2, $target2, "", "", $target2) | Out-Null

# --- Trailing space after the new hyperlink, matching the first one -------
$pLink2c = $d.Paragraphs.Item(25)
$pLink2c.Range.InsertAfter(" ")
